# Add a header row to the "基本資料" worksheet and fix the row-shift / parity bug.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("基本資料")

# Insert a new blank row above the existing row 1, shifting all data down
# (rows 1-4 become rows 2-5, keeping their original styles/values intact).
$ws.Rows("1:1").Insert()

# Populate the new header row (columns A-I only; column J stays blank).
$ws.Range("A1").Value = "Breed"
$ws.Range("B1").Value = "ID"
$ws.Range("C1").Value = "confusing_note"
$ws.Range("D1").Value = "Birthday"
$ws.Range("E1").Value = "Sire"
$ws.Range("F1").Value = "Dam"
$ws.Range("G1").Value = "reg_id"
$ws.Range("H1").Value = "Chinese_name"
$ws.Range("I1").Value = "Gender"
